# Matheus - Python.xlsx : translate headers/labels to pt-BR, tweak hours,
# move the total into F10, and drop the unused "Bom" (Good) cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- Header row (B1:E1) -----------------------------------------------
$ws.Range("B1").Value = "Tarefa"
$ws.Range("C1").Value = "Valor da Hora"
$ws.Range("D1").Value = "Horas Trabalhadas"
$ws.Range("E1").Value = "Valor"

# --- Task names (B2:B4) -------------------------------------------------
$ws.Range("B2").Value = "Parser Basketball ESPN"
$ws.Range("B3").Value = "Parser Baseball ESPN"
$ws.Range("B4").Value = "Parser Baseball CBS"

# --- Hours worked updates ------------------------------------------------
$ws.Range("D2").NumberFormat = "0"
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 2

# --- Move the grand total from E12 into F10 ------------------------------
$ws.Range("E12").Clear()
$ws.Range("F10").Formula = "=SUM(E2:E10)"
$ws.Range("F10").NumberFormat = "[$$-540A]#,##0.00"

# --- Remove the now-unused "Bom" (Good) cell style -----------------------
$wb.Styles("Bom").Delete()

# --- Selection / view bookkeeping ----------------------------------------
$ws.Range("B1").Select()
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("D4").Select()

$wb.RefreshAll()
